$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers in I1 and J1, matching the style used by the other
# header cells (e.g. H1) by copying its formatting over.
$ws.Range("H1:H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new columns I (I0) and J (IF), rows 2-43
$data = @(
    @(8, 8),
    @(9, 9),
    @(5, 5),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(6, 6),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(6, 6),
    @(6, 7),
    @(9, 9),
    @(8, 8),
    @(8, 9),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(7, 8),
    @(6, 7),
    @(9, 9),
    @(8, 8),
    @(6, 6),
    @(9, 9),
    @(9, 9),
    @(7, 7),
    @(7, 8),
    @(6, 6),
    @(5, 6),
    @(8, 9),
    @(7, 7),
    @(7, 7),
    @(3, 4),
    @(7, 8),
    @(9, 9),
    @(7, 7),
    @(7, 7),
    @(6, 7),
    @(6, 6),
    @(6, 6),
    @(6, 6)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
